$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$tcsN = $nm.Theme.ThemeColorScheme
$tcsN.Colors(3).RGB = 999999
